$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("A34").Value = "SAP Regression Automation"
$ws.Range("B34").Value = "Yes"
$ws.Range("C34").Value = "CO_NW_33"
$ws.Range("C35").Value = "CO_NW_34"
$ws.Range("D34").Value = 5411316
$ws.Range("E34").Value = "Validate the Order Status in DBS"
$ws.Range("E35").Value = "SAPValidations"

$ws.Range("A35").Value = "SAP Regression Automation"
$ws.Range("B35").Value = "Yes"
$ws.Range("D35").Value = 5411317
